$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 23: Exception log entry
$ws.Cells.Item(23, 1).Value = 23
$ws.Cells.Item(23, 2).Value = "Exception: Update image directory"
$ws.Cells.Item(23, 3).Value = "[Errno 13] Permission denied: 'c:\Users\user10\Desktop\Hobby\Programming\EEEY3 Project\Web App\Elephant_Web_App_v2\static/image uploads/end device 2/2021-09-15 16-18-18-x-whale - Copy.jpg'"
$ws.Cells.Item(23, 4).Value = "07/05/2022 03:36:31 AM"

# Row 24: Object Detection log entry
$ws.Cells.Item(24, 1).Value = 24
$ws.Cells.Item(24, 2).Value = "Object Detection"
$ws.Cells.Item(24, 3).Value = "New image from end device 2 detected and recorded to database."
$ws.Cells.Item(24, 4).Value = "07/05/2022 03:36:31 AM"
